$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account Advanced Find View")
$ws.Range("D2").Value = "Employer-WrongDataType"
$ws.Columns.Item(4).ColumnWidth = 23.1
$ws.Range("D3").Select() | Out-Null
